# Add two new columns, "I0" (col I) and "IF" (col J), to the header row
# plus the value 9 for each of the three data rows — mirrors the existing
# "IP" column (H) both in content and in header styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy H1's formatting (bold/border/centered header style) onto the new
# header cells, then overwrite their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Fill in the data rows (2-4) with the value 9 for both new columns.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
